# Add a new bullet/paragraph to the "Using RSTUDIO" slide's content
# placeholder: "There are 5000 observations and 7 variables"
#
# The slide is identified by its PowerPoint creationId (cId 2024379099 /
# sldId 272 in the source OOXML), which corresponds to slide index 6 in
# this deck ("Using RSTUDIO").

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(6)

$shape = $null
foreach ($candidate in $s.Shapes) {
    if ($candidate.Name -eq "Content Placeholder 2") {
        $shape = $candidate
        break
    }
}
if ($shape -eq $null) {
    $shape = $s.Shapes.Item(2)
}

$tr = $shape.TextFrame.TextRange

# Append a new paragraph after the existing text, matching the author's
# added bullet point.
[void]$tr.InsertAfter("`r" + "There are 5000 observations and 7 variables")
